$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab (workbook.xml sheet name changes from
# "UniformA-HW40.xpc" to "UniformA")
$ws.Name = "UniformA"

# Add row 16, continuing the pattern of rows 2-15: column A holds the
# zero-based index (styled like the rest of column A, re-using the same
# format as A15), column B repeats the last label
# ("HexGrid-60degTilt5degRes"), and columns C:P are all 1.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(16, 1).Value = 14

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
